$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Insert a new row above the current "Eslovaquia" row (row 63) to make
#    room for Serbia in its new, alphabetically-earlier position.
$ws.Rows("63").Insert()

# 2) Populate the newly inserted row with Serbia's updated statistics.
$ws.Range("A63").Value = "Serbia"
$ws.Range("B63").Value = 188
$ws.Range("C63").Value = 17
$ws.Range("D63").Value = 2
$ws.Range("E63").Value = 185
$ws.Range("F63").Value = 4
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 1

# 3) The old Serbia row has now shifted down to row 66 (Eslovaquia and
#    Kuwait moved down by one to rows 64 and 65). Remove that now-duplicate
#    Serbia row so the remaining countries (Bulgaria, etc.) return to their
#    original row numbers.
$ws.Rows("66").Delete()

# 4) Update the "last updated" timestamp shown in cell A1.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 08:16"
